$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.847.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.527.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.97%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.14%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.581'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.525.58'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.07%  '
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.982.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.806.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000142'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.520.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '334.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.51%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.32%  '
$ws.Range("B27").Value = 'SuiNetwork'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.32'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0808'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '177.39'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("E34").Value = '  +6.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '413.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.397'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.604'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0965'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("E49").Value = '  +5.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.76%  '
